$wb = $excel.ActiveWorkbook

# "Create Teams" sheet - replace the sample row with the real BU/team data
$wsCreate = $wb.Worksheets.Item("Create Teams")

$wsCreate.Range("A2").Value = "0-ES-BRJ-07"
$wsCreate.Range("B2").Value = "ETGA0001"
$wsCreate.Range("C2").Value = "ZPQ"
$wsCreate.Range("D2").Value = "T501"
$wsCreate.Range("E2").Value = "GAMESA"

# Activate "Create Teams" sheet and set the last selection the author left it in
$wsCreate.Activate()
$wsCreate.Range("B20").Select()

# "Assign Teams" sheet keeps its prior selection but is no longer the active tab
$wsAssign = $wb.Worksheets.Item("Assign Teams")
$wsAssign.Range("C9").Select()

$wsCreate.Activate()
